$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-9, columns G:T
# Row 2
$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5532856666666667
$ws.Range("N2").Value = 1.659857
$ws.Range("O2").Value = 0.4010144607159208
$ws.Range("P2").Value = 0.4010144607159208
$ws.Range("Q2").Value = 19.77379418586423
$ws.Range("R2").Value = 177.964147672778
$ws.Range("S2").Value = 0.007818696874606528
$ws.Range("T2").Value = 0.007818696874606528

# Row 3
$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8264293333333333
$ws.Range("N3").Value = 2.479288
$ws.Range("O3").Value = 0.5989855392840792
$ws.Range("P3").Value = 0.5989855392840792
$ws.Range("Q3").Value = 29.53563508150578
$ws.Range("R3").Value = 265.820715733552
$ws.Range("S3").Value = 0.01167859721460913
$ws.Range("T3").Value = 0.01167859721460913

# Row 4
$ws.Range("H4").Value = 5067.86792
$ws.Range("I4").Value = 0.9215900675332435
$ws.Range("J4").Value = 0.9215900675332435
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5532856666666667
$ws.Range("N4").Value = 1.659857
$ws.Range("O4").Value = 0.4010144607159208
$ws.Range("P4").Value = 0.4010144607159208
$ws.Range("Q4").Value = 934.6595602319377
$ws.Range("R4").Value = 8411.93604208744
$ws.Range("S4").Value = 0.3695709439329927
$ws.Range("T4").Value = 0.3695709439329927

# Row 5
$ws.Range("H5").Value = 5067.86792
$ws.Range("I5").Value = 0.9215900675332435
$ws.Range("J5").Value = 0.9215900675332435
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8264293333333333
$ws.Range("N5").Value = 2.479288
$ws.Range("O5").Value = 0.5989855392840792
$ws.Range("P5").Value = 0.5989855392840792
$ws.Range("Q5").Value = 1396.078235515662
$ws.Range("R5").Value = 12564.70411964096
$ws.Range("S5").Value = 0.5520191236002507
$ws.Range("T5").Value = 0.5520191236002507

# Row 6
$ws.Range("G6").Value = 93.641553
$ws.Range("H6").Value = 280.924659
$ws.Range("I6").Value = 0.05108605424341119
$ws.Range("J6").Value = 0.05108605424341119
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5532856666666667
$ws.Range("N6").Value = 1.659857
$ws.Range("O6").Value = 0.4010144607159208
$ws.Range("P6").Value = 0.4010144607159208
$ws.Range("Q6").Value = 51.810529079307
$ws.Range("R6").Value = 466.2947617137631
$ws.Range("S6").Value = 0.02048624649252581
$ws.Range("T6").Value = 0.02048624649252581

# Row 7
$ws.Range("G7").Value = 93.641553
$ws.Range("H7").Value = 280.924659
$ws.Range("I7").Value = 0.05108605424341119
$ws.Range("J7").Value = 0.05108605424341119
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8264293333333333
$ws.Range("N7").Value = 2.479288
$ws.Range("O7").Value = 0.5989855392840792
$ws.Range("P7").Value = 0.5989855392840792
$ws.Range("Q7").Value = 77.388126218088
$ws.Range("R7").Value = 696.4931359627921
$ws.Range("S7").Value = 0.03059980775088537
$ws.Range("T7").Value = 0.03059980775088537

# Row 8
$ws.Range("G8").Value = 14.34625366666667
$ws.Range("H8").Value = 43.038761
$ws.Range("I8").Value = 0.007826584134129748
$ws.Range("J8").Value = 0.007826584134129748
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5532856666666667
$ws.Range("N8").Value = 1.659857
$ws.Range("O8").Value = 0.4010144607159208
$ws.Range("P8").Value = 0.4010144607159208
$ws.Range("Q8").Value = 7.937576524130779
$ws.Range("R8").Value = 71.43818871717701
$ws.Range("S8").Value = 0.003138573415795823
$ws.Range("T8").Value = 0.003138573415795823

# Row 9
$ws.Range("G9").Value = 14.34625366666667
$ws.Range("H9").Value = 43.038761
$ws.Range("I9").Value = 0.007826584134129748
$ws.Range("J9").Value = 0.007826584134129748
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8264293333333333
$ws.Range("N9").Value = 2.479288
$ws.Range("O9").Value = 0.5989855392840792
$ws.Range("P9").Value = 0.5989855392840792
$ws.Range("Q9").Value = 11.85616485357422
$ws.Range("R9").Value = 106.705483682168
$ws.Range("S9").Value = 0.004688010718333925
$ws.Range("T9").Value = 0.004688010718333925
